# Update the two-digit multiplication problems to new values.
$d = $word.ActiveDocument

$replacements = @(
    @{ old = "72×71="; new = "23×96=" },
    @{ old = "31×96="; new = "41×67=" },
    @{ old = "87×76="; new = "46×12=" },
    @{ old = "11×37="; new = "30×41=" },
    @{ old = "67×26="; new = "52×41=" },
    @{ old = "40×41="; new = "57×75=" },
    @{ old = "74×51="; new = "72×91=" },
    @{ old = "92×23="; new = "90×27=" },
    @{ old = "40×25="; new = "89×40=" },
    @{ old = "78×55="; new = "54×48=" },
    @{ old = "78×47="; new = "26×76=" },
    @{ old = "38×60="; new = "58×50=" },
    @{ old = "33×26="; new = "42×81=" },
    @{ old = "93×56="; new = "79×78=" },
    @{ old = "16×34="; new = "87×57=" },
    @{ old = "53×17="; new = "21×88=" },
    @{ old = "75×48="; new = "55×91=" },
    @{ old = "22×56="; new = "39×71=" },
    @{ old = "15×58="; new = "20×13=" },
    @{ old = "70×80="; new = "65×84=" },
    @{ old = "38×20="; new = "81×90=" },
    @{ old = "77×92="; new = "54×80=" },
    @{ old = "98×45="; new = "29×36=" },
    @{ old = "11×69="; new = "34×79=" },
    @{ old = "52×72="; new = "67×48=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
